$d = $word.ActiveDocument

# Locate the paragraph containing the sentence to edit (rather than
# assuming a fixed paragraph index).
$p = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text -like "*Hawaii and Wyoming*") {
        $p = $candidate
        break
    }
}
if ($p -eq $null) {
    throw "Could not locate target paragraph"
}

# Target segments: same Times New Roman formatting throughout, but split
# into 5 separate runs (as produced by Word's own edit-tracking) because
# "Utah" -> "South Dakota" and "7.45" -> "8.59" were typed as separate
# edits inside the original sentence.
$segments = @(
    "By state, Hawaii and Wyoming had the greatest percentage drop in total annual payroll (11.48% and 11.25% respectively) while Delaware and ",
    "South Dakota",
    " had the highest percentage increase (9.74% and ",
    "8.59",
    "% respectively). "
)

$full = $p.Range.Text
$pStart = $p.Range.Start
$pEnd = $p.Range.Start + $full.Length

# Lay down all the text first (replacing the paragraph's whole range with
# segment 0, then inserting each subsequent segment right after the
# previous one), recording the [start,end) bounds of each segment.
$rng = $d.Range($pStart, $pEnd)
$rng.Text = $segments[0]

$bounds = @()
$s0 = $pStart
$e0 = $pStart + $segments[0].Length
$bounds += ,@($s0, $e0)

$cursor = $e0
for ($i = 1; $i -lt $segments.Length; $i++) {
    $insertRng = $d.Range($cursor, $cursor)
    $insertRng.Text = $segments[$i]
    $s = $cursor
    $e = $cursor + $segments[$i].Length
    $bounds += ,@($s, $e)
    $cursor = $e
}

# Toggle Bold on/off (net no-op formatting-wise) over each segment, walking
# from the LAST segment back to the FIRST. This forces the engine to keep
# each segment as its own <w:r> run instead of re-merging adjacent runs
# that share identical formatting.
for ($i = $bounds.Length - 1; $i -ge 0; $i--) {
    $b = $bounds[$i]
    $r = $d.Range($b[0], $b[1])
    $r.Bold = 1
    $r.Bold = 0
}

Write-Host "Final: [$($p.Range.Text)]"
